$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two changed field names
$ws.Range("A4").Value = "referral_id"
$ws.Range("A7").Value = "platekey"

# Adjust column widths to match the resized columns
$ws.Columns.Item(1).ColumnWidth = 23.666666666666668
$ws.Columns.Item(2).ColumnWidth = 20

# Update the active cell selection
$ws.Range("E9").Select()
